# Apply updated cryptocurrency price/volume data (and a few Coin/Link
# re-orderings) for the Sat Feb 11 18:43:41 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold numeric-/percent-looking text
# in this sheet (e.g. "308.10", "0.34%") -- force the Text number format
# before assigning so Excel keeps them as literal strings instead of
# converting them into numbers/percentages.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '308.10'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.34%'
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '40.83'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.16%'
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.121'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.43%'
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07612'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.01%'
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.282'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.78%'
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.617'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.22%'
# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.08%'
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9081'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.28%'
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1275'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '28.23%'
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1809'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.27%'
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09072'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.54%'
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04313'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-2.35%'
# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.86%'
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001253'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.44%'
# Row 16
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.04046'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.83%'
# Row 17
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.005829'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.48%'
# Row 18
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.351'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.58%'
# Row 19
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3314'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.48%'
# Row 20
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.970'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '3.35%'
# Row 21
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1394'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '4.13%'
# Row 22
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2708'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-4.74%'
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001273'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '4.85%'
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004033'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.90%'
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001272'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-2.10%'
# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '24.75%'
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02431'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '0.45%'
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05242'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '1.90%'
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007841'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.17%'
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1296'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.74%'
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006810'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-3.61%'
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001899'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-2.51%'
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007371'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-8.63%'
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3349'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.50%'
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006904'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '8.14%'
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000752'
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1074'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '1,595.13%'
